$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.311.57'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.85%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.440.12'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.98'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.02'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.72%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.443.49'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.33%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.556'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.30'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.89%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.426'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.49%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.037.64'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.33%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.34%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.38'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.56%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000175'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -6.84%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.373.66'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.421.19'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.15'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.68'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.45%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '379.34'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.84'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.70%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.56'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.72%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.520'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -5.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000117'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.59'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.57%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.11'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.48%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.68%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.31%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.02'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.17'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.53%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.60%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.05'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.864'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +11.78%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.83'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.825.43'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.64%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0732'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.18'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.70%  '

$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.64'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.81%  '

$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.99'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.64%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.46'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.44'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.93%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.51'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +10.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0308'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '336.85'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +6.19%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.86%  '

$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.34'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.76%  '

$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.103'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.65%  '
